$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 7, pushing existing rows down
$ws.Rows.Item(7).Resize(2).Insert()

$ws.Cells.Item(7, 1).Value = "Cấu hình Vlan Switch HPE1111"
$ws.Cells.Item(8, 1).Value = "Cấu hình Vlan Switch UPE5523"

$ws.Cells.Item(7, 2).Value = "abc"
$ws.Cells.Item(8, 2).Value = "cdf"

# Apply wrap text style to new B cells (matching style used for long answers)
$ws.Range("B7:B8").WrapText = $true

# Update selection to B13 (one past last data row) as shown in the diff
$ws.Range("B13").Select()
